$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B6 value
$ws.Range("B6").Value = 80219

# Clear C6 entirely (cell becomes empty, no shifting of other cells)
$ws.Range("C6").ClearContents()

# Update Q6 and R6 values (rounded coordinates)
$ws.Range("Q6").Value = 546617
$ws.Range("R6").Value = 6916220
